# Re-insert the "First Slide Insertion" slide as a new second slide.
#
# Target state (per the diff): the deck gains a brand-new slide (id=258)
# positioned between the existing first slide (id=256) and the existing
# "First Slide Insertion" slide (id=257, becomes the 3rd slide). The new
# slide re-creates that same "Title"/"Content Placeholder" slide (a
# title placeholder reading "First Slide Insertion" plus an empty
# generic content placeholder) using the "Title and Content" layout.

$p = $ppt.ActivePresentation

# Insert a brand-new slide at position 2 using the classic "Title, Content"
# autolayout (ppLayoutText = 2) so the placeholders come out as
# <p:ph type="title"/> and <p:ph idx="1"/>, matching the duplicated slide.
$newSlide = $p.Slides.Add(2, 2)

# Title placeholder -> same text as the existing "First Slide Insertion" slide.
$titleShape = $newSlide.Shapes.Item(1)
$titleShape.Name = "Title"
$titleShape.TextFrame.TextRange.Text = "First Slide Insertion"

# Generic content placeholder stays empty, just rename it to match.
$contentShape = $newSlide.Shapes.Item(2)
$contentShape.Name = "Content Placeholder"
